# Update code for T8
# Insert 5 new customer rows at the top of the customer list (right after the header row),
# shifting all existing rows down by 5, and populate the new rows with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows right after the header row (row 1); this pushes existing data rows
# (old rows 2..119) down to (new rows 7..124) automatically.
$ws.Rows("2:6").Insert()

# Row 2 - ngọc hân
$ws.Range("A2").Value2 = "KH"
$ws.Range("B2").Value2 = 395
$ws.Range("C2").Value2 = "ngọc hân"
$ws.Range("D2").Value2 = "SÓC TRĂNG"
$ws.Range("E2").Value2 = 0
$ws.Range("F2").Value2 = 0
$ws.Range("G2").Value2 = 0
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 35000000
$ws.Range("J2").Value2 = 8000000

# Row 3 - lan anh
$ws.Range("A3").Value2 = "KH"
$ws.Range("B3").Value2 = 394
$ws.Range("C3").Value2 = "lan anh"
$ws.Range("D3").Value2 = "SÓC TRĂNG"
$ws.Range("E3").Value2 = 0
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "0372781693"
$ws.Range("G3").Value2 = 0
$ws.Range("H3").Value2 = 0
$ws.Range("I3").Value2 = 2800000
$ws.Range("J3").Value2 = 0

# Row 4 - lâm thị duyên
$ws.Range("A4").Value2 = "KH"
$ws.Range("B4").Value2 = 393
$ws.Range("C4").Value2 = "lâm thị duyên"
$ws.Range("D4").Value2 = "SÓC TRĂNG"
$ws.Range("E4").Value2 = 0
$ws.Range("F4").Value2 = 0
$ws.Range("G4").Value2 = 0
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 500000
$ws.Range("J4").Value2 = 0

# Row 5 - mạch minh tâm
$ws.Range("A5").Value2 = "KH"
$ws.Range("B5").Value2 = 385
$ws.Range("C5").Value2 = "mạch minh tâm"
$ws.Range("D5").Value2 = "SÓC TRĂNG"
$ws.Range("E5").Value2 = 0
$ws.Range("F5").Value2 = 0
$ws.Range("G5").Value2 = 0
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 2300000
$ws.Range("J5").Value2 = 0

# Row 6 - trần lê ngọc bảo yến
$ws.Range("A6").Value2 = "KH"
$ws.Range("B6").Value2 = 384
$ws.Range("C6").Value2 = "trần lê ngọc bảo yến "
$ws.Range("D6").Value2 = "SÓC TRĂNG"
$ws.Range("E6").Value2 = 0
$ws.Range("F6").Value2 = 0
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 3000000
$ws.Range("J6").Value2 = 0
